$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Swap the header labels in D1 and E1
$ws.Range("D1").Value = "LESS THAN DEMAND PAYOUT%"
$ws.Range("E1").Value = "MORE THAN DEMAND PAYOUT%"

# Ensure the percent-looking text values stay as plain text, not auto-converted numbers
$pctRange = $ws.Range("D2:E6")
$pctRange.NumberFormat = "@"

# Row 2 (NEELU MALIK / LESS THAN DEMAND)
$ws.Range("C2").Value = 120186
$ws.Range("D2").Value = "1%"
$ws.Range("E2").Value = "0%"

# Row 3 (RUKHSAR KHAN / LESS THAN DEMAND)
$ws.Range("C3").Value = 333104
$ws.Range("D3").Value = "3%"
$ws.Range("E3").Value = "0%"

# Row 4 (VIJAY KHANNA / LESS THAN DEMAND)
$ws.Range("D4").Value = "1%"
$ws.Range("E4").Value = "0%"

# Row 5 (NEELU MALIK / MORE THAN DEMAND)
$ws.Range("C5").Value = 229950

# Row 6 (RUKHSAR KHAN / MORE THAN DEMAND)
$ws.Range("C6").Value = 411626
$ws.Range("E6").Value = "2%"
